# Updates the cryptocurrency price/volume table (cryptos list) to reflect
# the latest scrape: refreshed Price (col D) and Volume(1h) (col E) values,
# plus a handful of rows where the coin ranking shuffled (name/link/price/
# volume all change together): rows 32/33 (LidoDAOToken<->Monero),
# 41/42 (NEARProtocol<->VeChain) and 46/48/49 (EnergySwap->Algorand,
# Algorand->Aave, Aave->MultiversX).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-converted to numbers (losing the original text formatting).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "44.151.29"
$ws.Range("E2").Value = "  +1.56%  "

$ws.Range("D3").Value = "2.251.76"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "307.62"
$ws.Range("E5").Value = "  -4.41%  "

$ws.Range("D6").Value = "98.40"
$ws.Range("E6").Value = "  -2.36%  "

$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  -3.36%  "

$ws.Range("D10").Value = "35.52"
$ws.Range("E10").Value = "  -3.74%  "

$ws.Range("D11").Value = "0.0822"
$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("D12").Value = "7.31"
$ws.Range("E12").Value = "  -5.06%  "

$ws.Range("E13").Value = "  -1.77%  "

$ws.Range("D14").Value = "2.594.14"
$ws.Range("E14").Value = "  +0.49%  "

$ws.Range("D15").Value = "2.250.47"
$ws.Range("E15").Value = "  +0.47%  "

$ws.Range("D16").Value = "0.837"
$ws.Range("E16").Value = "  -1.80%  "

$ws.Range("D17").Value = "13.80"
$ws.Range("E17").Value = "  -2.09%  "

$ws.Range("D18").Value = "43.990.60"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("D19").Value = "13.03"
$ws.Range("E19").Value = "  -3.53%  "

$ws.Range("D20").Value = "0.0₃0974"
$ws.Range("E20").Value = "  -0.82%  "

$ws.Range("D21").Value = "6.31"
$ws.Range("E21").Value = "  -3.87%  "

$ws.Range("D22").Value = "65.37"
$ws.Range("E22").Value = "  +0.63%  "

$ws.Range("D23").Value = "242.65"
$ws.Range("E23").Value = "  +2.53%  "

$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  -6.80%  "

$ws.Range("D25").Value = "1.98"
$ws.Range("E25").Value = "  -7.97%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").Value = "10.10"
$ws.Range("E27").Value = "  +0.65%  "

$ws.Range("D28").Value = "2.14"
$ws.Range("E28").Value = "  -1.86%  "

$ws.Range("D29").Value = "36.49"
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").Value = "6.20"
$ws.Range("E30").Value = "  -1.47%  "

$ws.Range("D31").Value = "20.18"
$ws.Range("E31").Value = "  +0.16%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "156.47"
$ws.Range("E32").Value = "  -1.56%  "

$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").Value = "3.51"
$ws.Range("E33").Value = "  +13.66%  "

$ws.Range("D34").Value = "0.0824"
$ws.Range("E34").Value = "  -2.91%  "

$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("D37").Value = "0.107"
$ws.Range("E37").Value = "  -3.91%  "

$ws.Range("D38").Value = "1.84"
$ws.Range("E38").Value = "  -3.73%  "

$ws.Range("D39").Value = "15.46"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").Value = "3.87"
$ws.Range("E40").Value = "  -8.57%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0307"
$ws.Range("E41").Value = "  -3.09%  "

$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").Value = "3.36"
$ws.Range("E42").Value = "  -10.40%  "

$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").Value = "1.770.32"
$ws.Range("E44").Value = "  -1.28%  "

$ws.Range("D45").Value = "87.66"
$ws.Range("E45").Value = "  +6.73%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.193"
$ws.Range("E46").Value = "  -2.84%  "

$ws.Range("D47").Value = "5.15"
$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "101.29"
$ws.Range("E48").Value = "  -1.60%  "

$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "55.89"
$ws.Range("E49").Value = "  -4.15%  "

$ws.Range("D50").Value = "8.24"
$ws.Range("E50").Value = "  -1.88%  "

$ws.Range("D51").Value = "69.47"
$ws.Range("E51").Value = "  -6.74%  "

Write-Output "cryptos list updated"
